$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the split "27" "/" "09" "/" "2025" runs in the revision-history
#    table into a single run "27/09/2025".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("27/09/2025", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "27/09/2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Rewrite the "Introdução" body paragraph: new wording, and drop the bold /
#    28pt direct formatting that used to live on the paragraph mark (w:pPr/w:rPr)
#    and on the run itself.
# ---------------------------------------------------------------------------
$introOld = "Este documento consolida todos os requisitos do sistema de Loja Virtual, incluindo os não documentados nos casos de uso. Ele serve como base para o desenvolvimento e validação do sistema."
$introNew = "O levantamento de requisitos foi realizado inicialmente de forma individual por cada membro do grupo e, em seguida, consolidado em reunião colaborativa, na qual os requisitos foram discutidos, refinados e priorizados. Este documento reúne todos os requisitos do sistema de Loja Virtual, incluindo aqueles não contemplados nos casos de uso, servindo como base para o desenvolvimento e a validação do sistema."

$introPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Este documento consolida")) {
        $introPara = $cand
        break
    }
}

$introXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>' + $introNew + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$introPara.Range.InsertXML($introXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Move the <w:lastRenderedPageBreak/> marker: it currently sits on the run
#    "Alteração mensal de senha." and needs to move up onto the "Segurança"
#    heading run instead.
# ---------------------------------------------------------------------------
$segPara = $null
$sennhaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t.StartsWith("Segurança") -and $segPara -eq $null) {
        $segPara = $cand
    }
    if ($t.StartsWith("Alteração mensal de senha") -and $sennhaPara -eq $null) {
        $sennhaPara = $cand
    }
}

$segXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Segurança</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$segPara.Range.InsertXML($segXml) | Out-Null

$senhaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Alteração mensal de senha.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sennhaPara.Range.InsertXML($senhaXml) | Out-Null

Write-Output "done"
